# repull data, push all data, mean calculation
# Update the dSF (column F) values on Sheet1 to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = -2
    6  = 0
    12 = 1
    13 = -1
    14 = 1
    20 = 0
    29 = -1
    30 = -1
    33 = -1
    34 = -3
    38 = 1
    41 = -1
    42 = 2
    46 = 0
    48 = -4
    50 = -3
    52 = 0
    57 = 2
    60 = -2
    61 = -4
    62 = 4
    65 = -5
    66 = 0
    68 = 5
    69 = 8
    70 = 0
    71 = 2
    72 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
